$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.872.95"
$ws.Range("E2").Value = "  +3.34%  "

$ws.Range("D3").Value = "1.675.05"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'219.64"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  +2.17%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'29.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("D10").Value = "'0.0639"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.88%  "

$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").Value = "1.917.72"
$ws.Range("E12").Value = "  +2.92%  "

$ws.Range("D13").Value = "1.662.15"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.59%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.42%  "

$ws.Range("D16").Value = "'4.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.57%  "

$ws.Range("D17").Value = "30.842.57"
$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").Value = "'65.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "'246.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "

$ws.Range("E20").Value = "  +2.17%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("D23").Value = "'9.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "'158.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").Value = "'15.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.07%  "

$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("E28").Value = "  +1.43%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.11%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.24%  "

$ws.Range("D33").Value = "'3.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.75%  "

$ws.Range("D34").Value = "1.519.75"
$ws.Range("E34").Value = "  +6.50%  "

$ws.Range("E35").Value = "  +3.68%  "

$ws.Range("D36").Value = "'83.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.83%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "'0.606"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.09%  "

$ws.Range("E39").Value = "  +4.86%  "

$ws.Range("D40").Value = "'2.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("D41").Value = "'2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.24%  "

$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0503"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.85%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.836"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("E47").Value = "  +4.35%  "

$ws.Range("D48").Value = "'50.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.92%  "

$ws.Range("D49").Value = "1.811.79"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("D50").Value = "0.0₆0120"
$ws.Range("E50").Value = "  +8.37%  "

$ws.Range("D51").Value = "'93.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
